$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows right after the current row 9 (before current row 10),
# shifting all existing rows 10:107 down to 12:109. This is a weekly data refresh:
# two fresh price records are prepended, and the two oldest records fall off the
# bottom of the pre-existing range (they simply end up as the new last rows
# 108:109, since the whole block just shifts down by two rows).
$ws.Rows("10:11").Insert()

# --- New row 10 ---
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C10").Value = "Los Lagos"
$ws.Range("D10").Value = 44490
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100101
$ws.Range("H10").Value = "Berries"
$ws.Range("I10").Value = 100112025
$ws.Range("J10").Value = "Frutilla"
$ws.Range("K10").Value = "Sin especificar"
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 150
$ws.Range("N10").Value = 15000
$ws.Range("O10").Value = 15000
$ws.Range("P10").Value = 15000
$ws.Range("Q10").Value = "$/bandeja 7 kilos"
$ws.Range("R10").Value = "Provincia de Melipilla"
$ws.Range("S10").Value = 2143
$ws.Range("T10").Value = 7

# --- New row 11 ---
$ws.Range("A11").Value = 4
$ws.Range("B11").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C11").Value = "Los Lagos"
$ws.Range("D11").Value = 44490
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = "Fruta"
$ws.Range("G11").Value = 100101
$ws.Range("H11").Value = "Berries"
$ws.Range("I11").Value = 100112025
$ws.Range("J11").Value = "Frutilla"
$ws.Range("K11").Value = "Sin especificar"
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 300
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 13000
$ws.Range("P11").Value = 12500
$ws.Range("Q11").Value = "$/bandeja 7 kilos"
$ws.Range("R11").Value = "Provincia de Melipilla"
$ws.Range("S11").Value = 1786
$ws.Range("T11").Value = 7
